# Update Name of Algo
# This script updates the "C" column (algorithm result) values for a set of
# rows in the RandomForest imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -13.7245
$ws.Range("C12").Value = -11.31799999999999
$ws.Range("C15").Value = -13.3571
$ws.Range("C27").Value = -12.463
$ws.Range("C28").Value = -13.733
$ws.Range("C31").Value = -13.51779999999999
$ws.Range("C32").Value = -13.1187
$ws.Range("C36").Value = -11.86580000000001
$ws.Range("C38").Value = -11.91330000000001
$ws.Range("C46").Value = -14.55459999999999
$ws.Range("C54").Value = -13.1233
$ws.Range("C55").Value = -13.74260000000001
$ws.Range("C56").Value = -12.9349
$ws.Range("C67").Value = -11.599
$ws.Range("C69").Value = -11.95799999999999
$ws.Range("C72").Value = -11.7421
$ws.Range("C73").Value = -11.16890000000001
$ws.Range("C83").Value = -13.69750000000001
$ws.Range("C86").Value = -14.33299999999999
$ws.Range("C91").Value = -12.4953
$ws.Range("C93").Value = -10.4451
$ws.Range("C99").Value = -12.06410000000001
$ws.Range("C104").Value = -12.90930000000001
$ws.Range("C105").Value = -12.64340000000001
